$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: add hours worked (C19), description (D19), and release version (E19)
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = "Fixed  serous bug with app crashing."
$ws.Range("E19").Value = "1,71a"

# Update view: scroll position and active selection
$ws.Range("E20").Select()
